# Agregar nueva BT de PPR para terremoto
# Inserts a new nomenclature row (concepto / Factor PPR Terremoto / PPR /
# terremoto_factor_pprr) right above the existing "Gasto Total NT Terremoto"
# row in the V2 sheet, shifting every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V2")

# Insert a new row at 63 - it inherits formatting from the row above (row 62),
# which matches the target styling (s=29 / s=30 / s=30 / s=15 / s=29).
$ws.Rows.Item(63).Insert()

$ws.Cells.Item(63, 1).Value = "concepto"
$ws.Cells.Item(63, 2).Value = "Factor PPR Terremoto"
$ws.Cells.Item(63, 3).Value = "PPR"
$ws.Cells.Item(63, 4).Value = "terremoto_factor_pprr"

# Mirror the author's final selection on the sheet.
$ws.Activate()
$ws.Range("C50").Select()
